$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("other-models")
$ws.Activate()
Write-Host "Active sheet: " $ws.Name
